# Updates the Price (D) and Volume(1h) (E) columns for the cryptos list
# For D-column values, temporarily force text format so numeric-looking
# strings (e.g. "0.630", "1.00") are preserved exactly as text instead of
# being auto-coerced to a Number by Excel (which would drop trailing zeros).
# Style is reset to "Normal" right after so no stray per-cell style survives.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.348.52"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +5.94%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.252.01"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.10%  "

$ws.Range("E4").Value = "  +0.17%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.63"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.88%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.630"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.15%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "61.95"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.74%  "

$ws.Range("E8").Value = "  +0.10%  "

$ws.Range("E9").Value = "  +2.72%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "59.54"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.81%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0901"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.08%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.585.35"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.17%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.74"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.34%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "22.17"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.13%  "

$ws.Range("E16").Value = "  -1.74%  "

$ws.Range("E17").Value = "  +0.37%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.266.54"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.94%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "42.209.72"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.79%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0911"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.67%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.35"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.25%  "

$ws.Range("E22").Value = "  -0.49%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "251.94"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +8.84%  "

$ws.Range("E24").Value = "  +0.03%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.39"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.67%  "

$ws.Range("E26").Value = "  +1.44%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.72"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.51%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.145"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.31%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "169.06"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.30%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.12"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.03%  "

$ws.Range("E31").Value = "  -3.02%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.71"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.07%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.122"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.23%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.03"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.54%  "

$ws.Range("E35").Value = "  +2.32%  "

$ws.Range("E36").Value = "  +2.54%  "

$ws.Range("E37").Value = "  -4.87%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.73"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.58%  "

$ws.Range("E39").Value = "  -3.91%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.000267"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +37.85%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.32%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0241"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.01%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.87"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.89%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.58"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.71%  "

$ws.Range("E45").Value = "  +0.46%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "99.53"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.90%  "

$ws.Range("E47").Value = "  +5.00%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.481.30"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.49%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "16.56"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -7.69%  "

$ws.Range("E50").Value = "  +0.06%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "52.71"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.84%  "
